# Update cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.575.41"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "3.496.08"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "612.36"
$ws.Range("D6").Value = "189.49"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "0.213"
$ws.Range("E9").Value = "  -3.91%  "
$ws.Range("D10").Value = "0.648"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "53.02"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("E12").Value = "  -4.20%  "
$ws.Range("D13").Value = "9.48"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "4.054.22"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "599.24"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("D16").Value = "69.652.17"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.00"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "3.488.19"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").Value = "17.18"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("D23").Value = "105.42"
$ws.Range("E23").Value = "  +12.18%  "
$ws.Range("D24").Value = "5.12"
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("D25").Value = "4.68"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "3.03"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "6.95"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "4.15"
$ws.Range("E31").Value = "  +9.23%  "
$ws.Range("D32").Value = "12.49"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "63.43"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "3.68"
$ws.Range("E37").Value = "  +6.72%  "
$ws.Range("D38").Value = "3.623.73"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D42").Value = "0.0₃0773"
$ws.Range("E42").Value = "  -3.96%  "
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "2.91"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").Value = "8.78"
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "132.54"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("E51").Value = "  -7.99%  "

# Rows 40-41: Bittensor and InjectiveProtocol swapped rank position, with refreshed price/volume
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "36.78"
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "503.36"
$ws.Range("E41").Value = "  -7.04%  "
